$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45960,7645.37551201644,7055.16475752583,8828,2923.32153990868,47.9369290597712),
    @(45961,6991.14466106109,5847.45670502095,3620,3263.83259904437,228.803721002722),
    @(45962,2148.01081475917,3411.98179232531,9004,3263.72956273574,-97.0120268724565),
    @(45963,2126.45966694183,3597.71792817645,9004,3240.22935228605,-90.252196647396),
    @(45964,8027.73259147043,6851.51055981854,9004,4248.33305669586,87.3268173547666),
    @(45965,8517.17560294576,7765.96772974836,9004,4586.80238049549,139.532087926827),
    @(45966,8517.17560294576,8052.55755901719,9004,4586.80238049549,151.473330813028),
    @(45967,8525.26551938708,8075.37107298074,9004,4588.7587085231,152.50540756266),
    @(45968,8530.1901118937,7182.65622398135,9004,4619.38696702862,116.585132958749),
    @(45969,3118.30130226871,4914.56914645502,9004,4308.05385142552,9.10929157835577),
    @(45970,2984.26537892464,4866.96748994846,9004,4300.78774350411,6.82313472719064),
    @(45971,9278.26748203409,8142.62917404075,9004,5247.55363186726,182.757616912834),
    @(45972,9278.26878400516,8385.64577244628,9004,5247.55739597229,192.883465350774),
    @(45973,9278.26878400516,8434.30103105549,9004,5247.55739597229,194.910767792824)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
